$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh -- update Coin/Link/Price/Volume(1h) cells per latest snapshot.
# Columns D/E are plain-text cells (t="inlineStr" in the source). Any new Price value
# that LOOKS like a plain number (no thousands-dots, no extra formatting) would otherwise
# be auto-coerced to a numeric cell by Excel's normal typed-entry behaviour, silently
# dropping meaningful trailing zeros (e.g. "91.00" -> 91, "0.200" -> 0.2). Force those
# specific cells to Text format first so the literal string round-trips exactly.

$ws.Range("D2").Value = '93.851.18'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '3.470.14'
$ws.Range("E3").Value = '  +4.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.56'
$ws.Range("E5").Value = '  +3.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.81'
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("E7").Value = '  +7.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.397'
$ws.Range("E8").Value = '  +4.21%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.02'
$ws.Range("E10").Value = '  +10.63%  '
$ws.Range("D11").Value = '3.471.12'
$ws.Range("E11").Value = '  +4.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.36'
$ws.Range("E12").Value = '  +11.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.200'
$ws.Range("E13").Value = '  +3.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.34'
$ws.Range("E14").Value = '  +8.14%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.115.06'
$ws.Range("E15").Value = '  +4.17%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '93.612.68'
$ws.Range("E16").Value = '  +1.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000250'
$ws.Range("E17").Value = '  +3.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.42'
$ws.Range("E18").Value = '  +7.33%  '
$ws.Range("D19").Value = '3.457.58'
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.46'
$ws.Range("E20").Value = '  +11.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.82'
$ws.Range("E21").Value = '  +9.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.521'
$ws.Range("E22").Value = '  +17.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.42'
$ws.Range("E23").Value = '  +10.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '503.75'
$ws.Range("E24").Value = '  +3.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.85'
$ws.Range("E25").Value = '  +13.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000187'
$ws.Range("E26").Value = '  +2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.00'
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.26'
$ws.Range("E28").Value = '  +8.59%  '
$ws.Range("D29").Value = '3.630.34'
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.55'
$ws.Range("E30").Value = '  +4.24%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.143'
$ws.Range("E31").Value = '  +10.34%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.77'
$ws.Range("E33").Value = '  +5.63%  '
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.181'
$ws.Range("E34").Value = '  +7.81%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.989'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.562'
$ws.Range("E36").Value = '  +8.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.72'
$ws.Range("E37").Value = '  +5.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '572.88'
$ws.Range("E38").Value = '  +12.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.64'
$ws.Range("E39").Value = '  +4.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.43'
$ws.Range("E40").Value = '  +4.42%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.150'
$ws.Range("E42").Value = '  +3.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.923'
$ws.Range("E43").Value = '  +7.20%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0428'
$ws.Range("E44").Value = '  +10.67%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.75'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.72'
$ws.Range("E46").Value = '  +4.84%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.60'
$ws.Range("E47").Value = '  +3.99%  '
$ws.Range("B48").Value = 'MantraDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.63'
$ws.Range("E48").Value = '  +2.29%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.26'
$ws.Range("E49").Value = '  +6.18%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.14'
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.28'
$ws.Range("E51").Value = '  +3.08%  '
